$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells P1:R1, matching style of existing header cells ---
$ws.Range("A1").Copy()
$ws.Range("P1:R1").PasteSpecial(-4122)
$ws.Range("P1").Value = "projectimage"
$ws.Range("Q1").Value = "authoritytocreate"
$ws.Range("R1").Value = "Company_Name"

# --- Blank out P2:R10 (new empty columns for existing rows) ---
$ws.Range("P2:R10").Value = ""

# --- Row 11: new data row (mirrors row 2, but "sonu" variant) ---
$ws.Range("A11").Value = "ArcelorMittal SA sonu"
$ws.Range("B11").Value = "GASTBEL0009"
$ws.Range("C11").Value = "Flémalle"
$ws.Range("D11").Value = "Wallonie"
$ws.Range("E11").Value = "Belgium"
$ws.Range("F11").Value = "BEL"
$ws.Range("G11").Value = 56
$ws.Range("H11").Value = "Europe"
$ws.Range("I11").Value = "Western Europe"
$ws.Range("J11").Value = 50.594707
$ws.Range("K11").Value = 5.466776
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = 5000060650
$ws.Range("N11").Value = "ArcelorMittal Belgium SA sonu"
$ws.Range("O11").Value = 5000030093
$ws.Range("P11").Value = ""
$ws.Range("Q11").Value = ""
$ws.Range("R11").Value = ""

# --- Row 12: mostly empty row, J12/K12 hold a literal two-double-quote string ---
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = '""'
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = '""'
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("P12").Value = ""
$ws.Range("Q12").Value = ""
$ws.Range("R12").Value = ""

# --- Row 13: new data row (partial; several fields stored as text) ---
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = "Archiologist"
$ws.Range("C13").Value = "Dallas"
$ws.Range("D13").Value = "New York"
$ws.Range("E13").Value = "US"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "1234521"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "56"
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "40.92679582427576"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = "-94.53104228055014"
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = "Arcel"
$ws.Range("O13").Value = ""
$ws.Range("P13").Value = ""
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = ""
